# Update scripts with new TPM values.
# The underlying NATMI computation was re-run with updated TPM data, which
# changed the numeric results and also dropped the "MuSCs" target-cluster
# rows from the table (rows 8-10 in the old layout), leaving 6 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (old rows 8, 9, 10) - the recomputed
# table only keeps 6 data rows (plus the header).
$ws.Rows("8:10").Delete()

# Row 2 (ECs -> ECs): recalculated values
$ws.Range("G2").Value = 0.1079986666666667
$ws.Range("H2").Value = 0.323996
$ws.Range("I2").Value = 0.004187739561209694
$ws.Range("J2").Value = 0.004187739561209694
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05619066666666667
$ws.Range("N2").Value = 0.168572
$ws.Range("O2").Value = 0.3931387525216601
$ws.Range("P2").Value = 0.39313875252166
$ws.Range("Q2").Value = 0.006068517079111112
$ws.Range("R2").Value = 0.054616653712
$ws.Range("S2").Value = 0.001646362706979583
$ws.Range("T2").Value = 0.001646362706979583

# Row 3 (ECs -> FAPs): recalculated values
$ws.Range("G3").Value = 0.1079986666666667
$ws.Range("H3").Value = 0.323996
$ws.Range("I3").Value = 0.004187739561209694
$ws.Range("J3").Value = 0.004187739561209694
$ws.Range("O3").Value = 0.60686124747834
$ws.Range("P3").Value = 0.60686124747834
$ws.Range("Q3").Value = 0.009367552349777779
$ws.Range("R3").Value = 0.08430797114800001
$ws.Range("S3").Value = 0.002541376854230111
$ws.Range("T3").Value = 0.002541376854230111

# Row 4 (FAPs -> ECs): now replaces the old "ECs -> MuSCs" row
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 24.98233333333333
$ws.Range("H4").Value = 74.947
$ws.Range("I4").Value = 0.9687110856121154
$ws.Range("J4").Value = 0.9687110856121155
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05619066666666667
$ws.Range("N4").Value = 0.168572
$ws.Range("O4").Value = 0.3931387525216601
$ws.Range("P4").Value = 0.39313875252166
$ws.Range("Q4").Value = 1.403773964888889
$ws.Range("R4").Value = 12.633965684
$ws.Range("S4").Value = 0.3808378677514501
$ws.Range("T4").Value = 0.3808378677514501

# Row 5 (FAPs -> FAPs): recalculated values
$ws.Range("D5").Value = "FAPs"
$ws.Range("I5").Value = 0.9687110856121154
$ws.Range("J5").Value = 0.9687110856121155
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08673766666666667
$ws.Range("N5").Value = 0.260213
$ws.Range("O5").Value = 0.60686124747834
$ws.Range("P5").Value = 0.60686124747834
$ws.Range("Q5").Value = 2.166909301222222
$ws.Range("R5").Value = 19.502183711
$ws.Range("S5").Value = 0.5878732178606654
$ws.Range("T5").Value = 0.5878732178606655

# Row 6 (MuSCs -> ECs): now replaces the old "FAPs -> FAPs" row
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "ECs"
$ws.Range("G6").Value = 0.6989190000000001
$ws.Range("H6").Value = 2.096757
$ws.Range("I6").Value = 0.02710117482667488
$ws.Range("J6").Value = 0.02710117482667489
$ws.Range("M6").Value = 0.05619066666666667
$ws.Range("N6").Value = 0.168572
$ws.Range("O6").Value = 0.3931387525216601
$ws.Range("P6").Value = 0.39313875252166
$ws.Range("Q6").Value = 0.039272724556
$ws.Range("R6").Value = 0.353454521004
$ws.Range("S6").Value = 0.01065452206323038
$ws.Range("T6").Value = 0.01065452206323038

# Row 7 (MuSCs -> FAPs): now replaces the old "FAPs -> MuSCs" row
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("G7").Value = 0.6989190000000001
$ws.Range("H7").Value = 2.096757
$ws.Range("I7").Value = 0.02710117482667488
$ws.Range("J7").Value = 0.02710117482667489
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08673766666666667
$ws.Range("N7").Value = 0.260213
$ws.Range("O7").Value = 0.60686124747834
$ws.Range("P7").Value = 0.60686124747834
$ws.Range("Q7").Value = 0.06062260324900001
$ws.Range("R7").Value = 0.5456034292410001
$ws.Range("S7").Value = 0.0164466527634445
$ws.Range("T7").Value = 0.01644665276344451
